$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.519.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.761.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.52"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.46"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.761.64"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.50%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.49"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000276"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.60"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.400.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.766.17"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.577.90"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.62"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.720"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000147"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -8.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.34"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.90"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.913.65"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.66"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.60"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.23"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.14"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.730.82"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.81"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.58%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.137"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.85"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.312"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.75"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.94"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.93"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "399.42"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000269"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -9.26%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0353"
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.92"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.67%  "
